$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 122, shifting existing rows 122..158 down to 123..159.
$ws.Rows.Item(122).Insert()

# Populate the new row 122 with its data (constants shared across the table plus
# the row-specific values from the diff).
$ws.Cells.Item(122, 1).Value = 11
$ws.Cells.Item(122, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(122, 3).Value = "Bíobío"
$ws.Cells.Item(122, 4).Value = 44855
$ws.Cells.Item(122, 5).Value = 8
$ws.Cells.Item(122, 6).Value = 100112032
$ws.Cells.Item(122, 7).Value = "Zapallo italiano"
$ws.Cells.Item(122, 8).Value = "Sin especificar"
$ws.Cells.Item(122, 9).Value = "Primera"
$ws.Cells.Item(122, 10).Value = 100
$ws.Cells.Item(122, 11).Value = 14000
$ws.Cells.Item(122, 12).Value = 15000
$ws.Cells.Item(122, 13).Value = 14500
$ws.Cells.Item(122, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(122, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(122, 16).Value = 290
$ws.Cells.Item(122, 17).Value = 50
$ws.Cells.Item(122, 18).Value = "Hortaliza"
